# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in columns AC, AD, AE (row 1)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the formatting used by the other header cells (bold, bordered,
# centered/top-aligned) by copying the format from the neighboring
# header cell (AB1) onto the new header cells.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Fill in the team record for every player row (rows 2-42) with the
# 1991 Mets' record: 77 wins, 84 losses, 0 ties.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 29).Value = 77
    $ws.Cells.Item($r, 30).Value = 84
    $ws.Cells.Item($r, 31).Value = 0
}
